$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B2").Value = 21
$ws.Range("B3").Value = 2250000
$ws.Range("B13").Value = 1500000
$ws.Range("B23").Value = 2250000
$ws.Range("B31").Value = -2700000
$ws.Range("B32").Value = 1500000
$ws.Range("B33").Value = 2920000
$ws.Range("B34").Value = 1720000.000000001

$wb.Save()
